$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cs = $m.ColorScheme
$c1 = $cs.Colors(1)
$c1.RGB = 255
Write-Output "done"
